# Update the "Training Dashboard" sheet:
# - Decrease each "PERIOD TO EXPIRE" (column H) value by 1
# - Update "LAST UPDATE" (column I) date text from 03-Nov-2025 to 04-Nov-2025
# for rows 3 through 16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 16; $row++) {
    # Column H ("PERIOD TO EXPIRE") is a plain number - decrement it.
    # Use Value2 (not Value) since Value surfaces as a non-numeric Variant
    # wrapper here; Value2 reliably round-trips a real double.
    $hCell = $ws.Cells.Item($row, 8)
    $hCell.Value = $hCell.Value2 - 1

    # Column I ("LAST UPDATE") holds the date as literal text (not a real
    # date serial). Just assigning the "04-Nov-2025" string would get
    # auto-recognised as a date by the Value setter (like typing it into
    # Excel), flipping the cell to a numeric date serial with a new
    # date-formatted style. To keep it as plain text with the original
    # style untouched:
    #   1. Temporarily force the cell to Text format so entry isn't
    #      reinterpreted as a date.
    #   2. Write the literal text.
    #   3. Paste-format from an unaffected neighboring cell that already
    #      carries the original (unchanged) style, restoring it exactly.
    $iCell = $ws.Cells.Item($row, 9)
    $sameStyleCell = $ws.Cells.Item($row, 10)  # column J keeps its original style

    $iCell.NumberFormat = "@"
    $iCell.Value = "04-Nov-2025"

    $sameStyleCell.Copy()
    $iCell.PasteSpecial(-4122)  # xlPasteFormats
}

$excel.CutCopyMode = 0
